# This script updates the team-specific transition probability matrix
# (Markov chain row-stochastic matrix) on Sheet1 with refreshed values
# produced after simulating more games (see commit message: "added more
# games, sped up simulate game logic, and drafted optimization logic").
# Only the probability cells that changed between the previous and new
# simulation runs are touched; zero/unchanged cells are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("B2").Value = 0.1968599033816425
    $ws.Range("C2").Value = 0.5628019323671497
    $ws.Range("J2").Value = 0.01570048309178744
    $ws.Range("P2").Value = 0.1461352657004831
    $ws.Range("S2").Value = 0.0785024154589372
    $ws.Range("B3").Value = 0.008368200836820083
    $ws.Range("C3").Value = 0.03347280334728033
    $ws.Range("J3").Value = 0.01673640167364017
    $ws.Range("P3").Value = 0.7594142259414226
    $ws.Range("S3").Value = 0.1820083682008368
    $ws.Range("J4").Value = 0.0392156862745098
    $ws.Range("P4").Value = 0.5882352941176471
    $ws.Range("S4").Value = 0.3725490196078431
    $ws.Range("P5").Value = 0.3333333333333333
    $ws.Range("S5").Value = 0.6666666666666666
    $ws.Range("B6").Value = 0.07023411371237458
    $ws.Range("D6").Value = 0.01337792642140468
    $ws.Range("E6").Value = 0.001672240802675585
    $ws.Range("F6").Value = 0.07859531772575251
    $ws.Range("J6").Value = 0.2508361204013378
    $ws.Range("O6").Value = 0.02006688963210702
    $ws.Range("Q6").Value = 0.1220735785953177
    $ws.Range("R6").Value = 0.08193979933110368
    $ws.Range("S6").Value = 0.3612040133779264
    $ws.Range("B7").Value = 0.1038062283737024
    $ws.Range("D7").Value = 0.01730103806228374
    $ws.Range("F7").Value = 0.06920415224913495
    $ws.Range("J7").Value = 0.1262975778546713
    $ws.Range("O7").Value = 0.02249134948096886
    $ws.Range("Q7").Value = 0.171280276816609
    $ws.Range("R7").Value = 0.08823529411764706
    $ws.Range("S7").Value = 0.4013840830449827
    $ws.Range("B8").Value = 0.1003110419906687
    $ws.Range("D8").Value = 0.02021772939346812
    $ws.Range("F8").Value = 0.0583203732503888
    $ws.Range("J8").Value = 0.1088646967340591
    $ws.Range("O8").Value = 0.02721617418351478
    $ws.Range("Q8").Value = 0.1594090202177294
    $ws.Range("R8").Value = 0.08942457231726283
    $ws.Range("S8").Value = 0.4362363919129083
    $ws.Range("B9").Value = 0.09386973180076628
    $ws.Range("D9").Value = 0.01149425287356322
    $ws.Range("E9").Value = 0.001915708812260536
    $ws.Range("F9").Value = 0.06704980842911877
    $ws.Range("J9").Value = 0.1360153256704981
    $ws.Range("O9").Value = 0.0210727969348659
    $ws.Range("Q9").Value = 0.1360153256704981
    $ws.Range("R9").Value = 0.09003831417624521
    $ws.Range("S9").Value = 0.4425287356321839
    $ws.Range("B10").Value = 0.107492795389049
    $ws.Range("D10").Value = 0.01642651296829971
    $ws.Range("E10").Value = 0.0005763688760806917
    $ws.Range("F10").Value = 0.06685878962536022
    $ws.Range("J10").Value = 0.1216138328530259
    $ws.Range("O10").Value = 0.01613832853025936
    $ws.Range("Q10").Value = 0.2175792507204611
    $ws.Range("R10").Value = 0.07694524495677234
    $ws.Range("S10").Value = 0.3763688760806916
    $ws.Range("G11").Value = 0.1498172959805116
    $ws.Range("J11").Value = 0.06333739342265529
    $ws.Range("K11").Value = 0.1961023142509135
    $ws.Range("L11").Value = 0.5749086479902558
    $ws.Range("S11").Value = 0.01583434835566382
    $ws.Range("G12").Value = 0.75
    $ws.Range("J12").Value = 0.1788617886178862
    $ws.Range("K12").Value = 0.006097560975609756
    $ws.Range("L12").Value = 0.03252032520325204
    $ws.Range("S12").Value = 0.03252032520325204
    $ws.Range("G13").Value = 0.7153284671532847
    $ws.Range("J13").Value = 0.2627737226277372
    $ws.Range("S13").Value = 0.0218978102189781
    $ws.Range("F15").Value = 0.01006711409395973
    $ws.Range("H15").Value = 0.1644295302013423
    $ws.Range("I15").Value = 0.04865771812080537
    $ws.Range("J15").Value = 0.3557046979865772
    $ws.Range("K15").Value = 0.05704697986577181
    $ws.Range("M15").Value = 0.01342281879194631
    $ws.Range("N15").Value = 0.001677852348993289
    $ws.Range("O15").Value = 0.06711409395973154
    $ws.Range("S15").Value = 0.2818791946308725
    $ws.Range("F16").Value = 0.01518026565464896
    $ws.Range("H16").Value = 0.1840607210626186
    $ws.Range("I16").Value = 0.07969639468690702
    $ws.Range("J16").Value = 0.4098671726755218
    $ws.Range("K16").Value = 0.1290322580645161
    $ws.Range("M16").Value = 0.01707779886148008
    $ws.Range("N16").Value = 0.00189753320683112
    $ws.Range("O16").Value = 0.04174573055028463
    $ws.Range("S16").Value = 0.1214421252371917
    $ws.Range("F17").Value = 0.01588628762541806
    $ws.Range("H17").Value = 0.1831103678929766
    $ws.Range("I17").Value = 0.09113712374581939
    $ws.Range("J17").Value = 0.4255852842809364
    $ws.Range("K17").Value = 0.09949832775919733
    $ws.Range("M17").Value = 0.02675585284280936
    $ws.Range("N17").Value = 0.002508361204013378
    $ws.Range("O17").Value = 0.06103678929765886
    $ws.Range("S17").Value = 0.09448160535117058
    $ws.Range("F18").Value = 0.02103250478011472
    $ws.Range("H18").Value = 0.1759082217973231
    $ws.Range("I18").Value = 0.08795411089866156
    $ws.Range("J18").Value = 0.4225621414913958
    $ws.Range("K18").Value = 0.1108986615678776
    $ws.Range("M18").Value = 0.01529636711281071
    $ws.Range("O18").Value = 0.06500956022944551
    $ws.Range("S18").Value = 0.1013384321223709
    $ws.Range("F19").Value = 0.01658564483843294
    $ws.Range("H19").Value = 0.2230483271375465
    $ws.Range("I19").Value = 0.08607377752359165
    $ws.Range("J19").Value = 0.372319130683443
    $ws.Range("K19").Value = 0.1078066914498141
    $ws.Range("M19").Value = 0.02287675150128682
    $ws.Range("N19").Value = 0.001143837575064341
    $ws.Range("O19").Value = 0.06348298541607092
    $ws.Range("S19").Value = 0.1066628538747498

